# Updates the cryptos list data (Price / Volume(1h) columns, plus a few
# coin-rank position swaps) to match the latest scrape.
#
# Note: Price-column values are always text in this sheet (e.g. "68.577.18",
# "603.83"), even when they look like plain numbers. Excel's COM layer
# auto-coerces bare numeric-looking strings (e.g. "601.07") into real
# numbers when assigned via .Value, so for those we prefix with a leading
# apostrophe (the standard "force text" convention) before assigning.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Price {
    param(
        [int]$Row,
        [string]$Price
    )
    if ($Price -match '^[+-]?\d+(\.\d+)?$') {
        $ws.Range("D$Row").Value = "'" + $Price
    } else {
        $ws.Range("D$Row").Value = $Price
    }
}

function Set-Row {
    param(
        [int]$Row,
        [string]$Coin,
        [string]$Link,
        [string]$Price,
        [string]$Volume
    )
    if ($Coin -ne "") { $ws.Range("B$Row").Value = $Coin }
    if ($Link -ne "") { $ws.Range("C$Row").Value = $Link }
    if ($Price -ne "") { Set-Price $Row $Price }
    if ($Volume -ne "") { $ws.Range("E$Row").Value = $Volume }
}

# Row 2 - Bitcoin
Set-Row 2 "" "" "68.571.39" "  -0.82%  "
# Row 3 - Ethereum
Set-Row 3 "" "" "3.912.99" "  +4.06%  "
# Row 4 - TetherUSD
Set-Row 4 "" "" "" "  -0.02%  "
# Row 5 - BNB
Set-Row 5 "" "" "601.07" "  -0.29%  "
# Row 6 - Solana
Set-Row 6 "" "" "165.34" "  -0.93%  "
# Row 7 - LidoStakedEther
Set-Row 7 "" "" "3.912.06" "  +4.09%  "
# Row 8 - USDC
Set-Row 8 "" "" "" "  +0.00%  "
# Row 9 - XRP
Set-Row 9 "" "" "" "  -1.92%  "
# Row 10 - Dogecoin
Set-Row 10 "" "" "" "  -3.57%  "
# Row 11 - Toncoin
Set-Row 11 "" "" "6.39" "  +0.03%  "
# Row 12 - Cardano
Set-Row 12 "" "" "" "  -0.12%  "
# Row 13 - Avalanche
Set-Row 13 "" "" "36.93" "  -2.18%  "
# Row 14 - ShibaInu
Set-Row 14 "" "" "" "  -1.06%  "
# Row 15 - WrappedliquidstakedEther2.0
Set-Row 15 "" "" "4.565.82" "  +4.02%  "
# Row 16 - WrappedEther
Set-Row 16 "" "" "3.933.23" "  +4.06%  "
# Row 17 - WrappedBTC
Set-Row 17 "" "" "68.783.06" "  -0.54%  "
# Row 18 - Polkadot
Set-Row 18 "" "" "7.41" "  +0.11%  "
# Row 19 - TRON
Set-Row 19 "" "" "" "  -0.82%  "
# Row 20 - Chainlink
Set-Row 20 "" "" "16.98" "  -4.38%  "
# Row 21 - Uniswap
Set-Row 21 "" "" "11.13" "  -1.73%  "
# Row 22 - BitcoinCash
Set-Row 22 "" "" "484.16" "  -1.53%  "

# Rows 23/24 swap: Polygon and PEPE swap rank positions
Set-Row 23 "PEPE" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe" "0.0000170" "  +14.10%  "
Set-Row 24 "Polygon" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic" "0.718" "  -1.32%  "

# Row 25 - Litecoin
Set-Row 25 "" "" "84.38" "  -0.38%  "
# Row 26 - Fetch.AI
Set-Row 26 "" "" "2.24" "  -1.51%  "
# Row 27 - InternetComputer(DFINITY)
Set-Row 27 "" "" "12.01" "  -2.09%  "
# Row 28 - RenderToken
Set-Row 28 "" "" "10.11" "  +0.51%  "
# Row 29 - Dai
Set-Row 29 "" "" "" "  -0.06%  "
# Row 30 - PancakeSwap
Set-Row 30 "" "" "" "  -1.04%  "
# Row 31 - WrappedeETH
Set-Row 31 "" "" "4.061.21" "  +4.01%  "
# Row 32 - NEARProtocol
Set-Row 32 "" "" "7.85" "  -3.34%  "
# Row 33 - ImmutableX
Set-Row 33 "" "" "2.37" "  -2.34%  "
# Row 34 - EthereumClassic
Set-Row 34 "" "" "31.92" "  +0.51%  "
# Row 35 - RenzoRestakedETH
Set-Row 35 "" "" "3.854.57" "  +4.17%  "
# Row 36 - Hedera
Set-Row 36 "" "" "" "  -1.01%  "
# Row 37 - Mantle (unchanged)
# Row 38 - Kaspa
Set-Row 38 "" "" "0.140" "  +0.00%  "
# Row 39 - Filecoin
Set-Row 39 "" "" "5.86" "  -1.38%  "
# Row 40 - FirstDigitalUSD
Set-Row 40 "" "" "" "  -0.04%  "

# Rows 41/42/43 rotate: TheGraph, Bittensor, dogwifhat -> dogwifhat, TheGraph, Bittensor
Set-Row 41 "dogwifhat" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif" "3.04" "  -2.39%  "
Set-Row 42 "TheGraph" "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt" "0.319" "  -2.40%  "
Set-Row 43 "Bittensor" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao" "435.38" "  +1.80%  "

# Row 44 - OKB
Set-Row 44 "" "" "48.49" "  -0.25%  "
# Row 45 - Stacks
Set-Row 45 "" "" "" "  -0.88%  "
# Row 46 - USDe
Set-Row 46 "" "" "" "  +0.00%  "
# Row 47 - Cosmos
Set-Row 47 "" "" "8.45" "  +0.16%  "
# Row 48 - Maker
Set-Row 48 "" "" "2.830.80" "  +1.10%  "
# Row 49 - Monero
Set-Row 49 "" "" "141.94" "  -0.51%  "
# Row 50 - EnergySwap
Set-Row 50 "" "" "25.93" "  +8.99%  "
# Row 51 - VeChain
Set-Row 51 "" "" "0.0353" "  -0.05%  "
